# Fixed the dialog window in PrevSO() - added uneditable fields
#
# The SO dialog now treats DEPOSIT/PRICE as read-only text fields instead of
# editable numerics, and a couple of customer-name entries got trimmed while
# testing the fix. Reproduce the resulting cell values:
#   C5: "Ashley Baker" -> "Ashley Bake"
#   G5: 55    (number) -> "55.0" (text)
#   H5: 120.99(number) -> "120.99" (text)
#   C9: "Ashley" -> "AshleyB"
#   G9: 15    (number) -> "15.0" (text)
#   H9: 74    (number) -> "74.0" (text)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $text) {
    # Assigning a numeric-looking string to .Value lets Excel "smart" parse
    # it back into a number. Temporarily switching the cell to the Text
    # number format forces the literal string to stick, then the original
    # number format is restored so the cell's style/appearance is unchanged.
    $originalFormat = $range.NumberFormat()
    $range.NumberFormat = "@"
    $range.Value = $text
    $range.NumberFormat = $originalFormat
}

# Row 5 (SO240111001 - Ashley Baker / Siamese Dream)
$ws.Range("C5").Value = "Ashley Bake"
Set-TextValue $ws.Range("G5") "55.0"
Set-TextValue $ws.Range("H5") "120.99"

# Row 9 (SO240112003 - Ashley / Chicago - EH)
$ws.Range("C9").Value = "AshleyB"
Set-TextValue $ws.Range("G9") "15.0"
Set-TextValue $ws.Range("H9") "74.0"
